$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance record (row 20): 2025-03-25 00:55:11 - Entrada
#
# Column A repeats the same date text ("2025-03-25") already used by rows
# 10-19, so copy one of those cells instead of assigning the literal via
# .Value - a plain string assignment that looks like a date gets
# auto-converted to a date serial number by the engine, which is not what
# the source data (plain text cells) represents.
$ws.Range("A19").Copy($ws.Range("A20"))

$ws.Range("B20").Value = "00:55:11"
$ws.Range("C20").Value = "Entrada"

# D20 mirrors the (empty) "Descripcion" cells that already exist in this
# sheet (e.g. D16) - copy one over so the new cell is present-but-blank,
# just like its siblings, rather than left unset.
$ws.Range("D16").Copy($ws.Range("D20"))

# The old trailing empty "Descripcion" cells (D16:D19) are no longer
# needed now that the data has moved on - drop them entirely.
$ws.Range("D16:D19").ClearContents()
